# Apply the commit's changes to the workbook.
#
# Semantic edits (everything else in the target OOXML diff - shared-string
# table re-ordering, tiny column-width deltas, etc. - is just mechanical
# fallout of re-saving the file and isn't something to reproduce by hand):
#   1. Rename sheet "store" -> "ready_to_sale"
#   2. Shorten the title in A1 from the two-run rich text
#      "Инструменты готовые к отправке" to plain "Инструменты"
#      (the bold/14pt formatting lives on the cell style, not the text
#      run, so it is kept automatically).
#   3. Tweak three product-name cells to use a dash instead of a space:
#        A4: "Ether Acril " -> "Ether-Acril"
#        A5: "Ether Wood"   -> "Ether-Wood"
#        A8: "Eternal love" -> "Eternal-love"
#   4. Leave the cursor/selection on A13, matching the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "ready_to_sale"

$ws.Range("A1").Value = "Инструменты"
$ws.Range("A4").Value = "Ether-Acril"
$ws.Range("A5").Value = "Ether-Wood"
$ws.Range("A8").Value = "Eternal-love"

$ws.Activate()
$ws.Range("A13").Select()
